$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.979.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "'1.716.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.26%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'317.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "'0.3974"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("D8").Value = "'0.4118"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.33%  "

$ws.Range("D9").Value = "'1.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").Value = "'1.005"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.46%  "

$ws.Range("D11").Value = "'53.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.29%  "

$ws.Range("D12").Value = "'0.08947"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.90%  "

$ws.Range("D13").Value = "'7.715"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.91%  "

$ws.Range("D14").Value = "'24.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.34%  "

$ws.Range("D15").Value = "'8.163"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "'0.00001373"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.06%  "

$ws.Range("D17").Value = "'1.686.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").Value = "'100.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").Value = "'0.07150"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("D20").Value = "'20.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.19%  "

$ws.Range("D21").Value = "'7.506"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.83%  "

$ws.Range("D22").Value = "'1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").Value = "'14.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "

$ws.Range("D24").Value = "'24.977.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("D25").Value = "'3.147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").Value = "'2.332"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "'23.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "

$ws.Range("D28").Value = "'9.299"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +24.05%  "

$ws.Range("D29").Value = "'165.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.23%  "

$ws.Range("D30").Value = "'140.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("D31").Value = "'5.231"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("D32").Value = "'7.899"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.15%  "

$ws.Range("D33").Value = "'0.09036"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.76%  "

$ws.Range("D34").Value = "'1.873.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("D35").Value = "'1.088"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").Value = "'0.03013"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.73%  "

$ws.Range("D37").Value = "'0.2811"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.00%  "

$ws.Range("E38").Value = "  -4.34%  "

$ws.Range("D39").Value = "'1.970"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.36%  "

$ws.Range("D40").Value = "'14.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.90%  "

$ws.Range("D41").Value = "'0.09327"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("D42").Value = "'0.8132"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.99%  "

$ws.Range("D43").Value = "'1.492"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.92%  "

$ws.Range("D44").Value = "'16.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.54%  "

$ws.Range("D45").Value = "'0.7391"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").Value = "'2.651"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").Value = "'4.268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("D48").Value = "'1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").Value = "'1.351"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("D50").Value = "'140.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").Value = "'93.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.41%  "
